# chore: update Sheets via scheduled runner
# Refreshes cached market-board price/profit figures (columns H-N:
# currentAveragePrice(NQ/HQ), LevePrice(NQ/HQ), LeveProfit(NQ/HQ)) across
# several crafting-leve worksheets to the latest scraped snapshot.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1615.2307
$ws.Range("I40").Value = 1640
$ws.Range("J40").Value = 1532.6666
$ws.Range("K40").Value = 1640
$ws.Range("L40").Value = 1532.6666
$ws.Range("M40").Value = -1465
$ws.Range("N40").Value = -1882.6666

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 587
$ws.Range("I5").Value = 716
$ws.Range("K5").Value = 716
$ws.Range("M5").Value = -604
$ws.Range("H38").Value = 16804
$ws.Range("I38").Value = 5340
$ws.Range("J38").Value = 34000
$ws.Range("K38").Value = 5340
$ws.Range("L38").Value = 34000
$ws.Range("M38").Value = -4873
$ws.Range("N38").Value = -34934
$ws.Range("H74").Value = 1681.1111
$ws.Range("I74").Value = 1532.2858
$ws.Range("J74").Value = 2202
$ws.Range("K74").Value = 1532.2858
$ws.Range("L74").Value = 2202
$ws.Range("M74").Value = -658.2858000000001
$ws.Range("N74").Value = -3950
$ws.Range("H77").Value = 1681.1111
$ws.Range("I77").Value = 1532.2858
$ws.Range("J77").Value = 2202
$ws.Range("K77").Value = 7661.429
$ws.Range("L77").Value = 11010
$ws.Range("M77").Value = -3293.429
$ws.Range("N77").Value = -19746
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 2271.7273
$ws.Range("I132").Value = 2034.8572
$ws.Range("J132").Value = 2686.25
$ws.Range("K132").Value = 6104.571599999999
$ws.Range("L132").Value = 8058.75
$ws.Range("M132").Value = -3574.571599999999
$ws.Range("N132").Value = -13118.75

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 587
$ws.Range("I4").Value = 716
$ws.Range("K4").Value = 716
$ws.Range("M4").Value = -601

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1400
$ws.Range("I31").Value = 1400
$ws.Range("K31").Value = 1400
$ws.Range("M31").Value = -1105
$ws.Range("H34").Value = 1400
$ws.Range("I34").Value = 1400
$ws.Range("K34").Value = 1400
$ws.Range("M34").Value = -1198
$ws.Range("H43").Value = 27666.111
$ws.Range("J43").Value = 27666.111
$ws.Range("L43").Value = 27666.111
$ws.Range("N43").Value = -28034.111
$ws.Range("H58").Value = 11499.5
$ws.Range("I58").Value = 9000
$ws.Range("K58").Value = 9000
$ws.Range("M58").Value = -8797
$ws.Range("H74").Value = 89314
$ws.Range("J74").Value = 89314
$ws.Range("L74").Value = 89314
$ws.Range("N74").Value = -91062
$ws.Range("H77").Value = 89314
$ws.Range("J77").Value = 89314
$ws.Range("L77").Value = 267942
$ws.Range("N77").Value = -276678
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("H95").Value = 35314.715
$ws.Range("J95").Value = 35314.715
$ws.Range("L95").Value = 35314.715
$ws.Range("N95").Value = -40806.715
$ws.Range("H99").Value = 2003235
$ws.Range("I99").Value = 1671666.6
$ws.Range("J99").Value = 2500587.5
$ws.Range("K99").Value = 1671666.6
$ws.Range("L99").Value = 2500587.5
$ws.Range("M99").Value = -1670168.6
$ws.Range("N99").Value = -2503583.5
$ws.Range("H101").Value = 27666.111
$ws.Range("J101").Value = 27666.111
$ws.Range("L101").Value = 27666.111
$ws.Range("N101").Value = -34156.111
$ws.Range("H107").Value = 949
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 949
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 949
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -4789
$ws.Range("H126").Value = 2003235
$ws.Range("I126").Value = 1671666.6
$ws.Range("J126").Value = 2500587.5
$ws.Range("K126").Value = 5014999.800000001
$ws.Range("L126").Value = 7501762.5
$ws.Range("M126").Value = -5012529.800000001
$ws.Range("N126").Value = -7506702.5
$ws.Range("H132").Value = 1718.6875
$ws.Range("I132").Value = 1718.6875
$ws.Range("K132").Value = 5156.0625
$ws.Range("M132").Value = -2626.0625
$ws.Range("H136").Value = 11499.5
$ws.Range("I136").Value = 9000
$ws.Range("K136").Value = 27000
$ws.Range("M136").Value = -24450

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 1080
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H60").Value = 153
$ws.Range("J60").Value = 200
$ws.Range("L60").Value = 600
$ws.Range("N60").Value = -1102
$ws.Range("H87").Value = 1014
$ws.Range("I87").Value = 1014
$ws.Range("K87").Value = 3042
$ws.Range("M87").Value = -1794
$ws.Range("H90").Value = 1014
$ws.Range("I90").Value = 1014
$ws.Range("K90").Value = 9126
$ws.Range("M90").Value = -2886

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 58333.332
$ws.Range("J93").Value = 58333.332
$ws.Range("L93").Value = 58333.332
$ws.Range("N93").Value = -62077.332
$ws.Range("H104").Value = 40000
$ws.Range("J104").Value = 40000
$ws.Range("L104").Value = 40000
$ws.Range("N104").Value = -46988
$ws.Range("H107").Value = 1815.25
$ws.Range("I107").Value = 631
$ws.Range("K107").Value = 631
$ws.Range("M107").Value = 1289

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").ClearContents()
$ws.Range("H24").Value = 19230.77
$ws.Range("I24").Value = 19230.77
$ws.Range("M24").Value = -19000.77
$ws.Range("H81").Value = 3759.4
$ws.Range("I81").Value = 3999.25
$ws.Range("J81").Value = 2800
$ws.Range("K81").Value = 7998.5
$ws.Range("L81").Value = 5600
$ws.Range("M81").Value = -6937.5
$ws.Range("N81").Value = -7722
$ws.Range("H84").Value = 3759.4
$ws.Range("I84").Value = 3999.25
$ws.Range("J84").Value = 2800
$ws.Range("K84").Value = 39992.5
$ws.Range("L84").Value = 28000
$ws.Range("M84").Value = -34688.5
$ws.Range("N84").Value = -38608
